$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Mobile_Phone value gets an extra number-only variant appended
$ws.Range("Q2").Value = "410-5644639, 4105644639"

# Row 3: Address_Line_2 becomes a bare number, Work_Phone becomes digits-only number
$ws.Range("K3").Value = 4
$ws.Range("P3").Value = 4074440909

# Row 4: last_name cleared entirely, Mobile_Phone becomes digits-only number
$ws.Range("E4").Clear()
$ws.Range("Q4").Value = 4077217359

# Row 4 height bumped slightly (matches the re-autofit after the edits)
$ws.Rows.Item(4).RowHeight = 24

# Selection moves to Z5
$ws.Range("Z5").Select()
